$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - column F holds "想去人数" (want-to-go count)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 23
$wsExpo.Range("F5").Value = 300
$wsExpo.Range("F7").Value = 1039
$wsExpo.Range("F9").Value = 536
$wsExpo.Range("F11").Value = 166
$wsExpo.Range("F12").Value = 13271
$wsExpo.Range("F16").Value = 5462
$wsExpo.Range("F17").Value = 5565
$wsExpo.Range("F18").Value = 35

# Sheet "全部类型" (all types) - same events duplicated with offset rows
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F10").Value = 23
$wsAll.Range("F21").Value = 300
$wsAll.Range("F29").Value = 1039
$wsAll.Range("F31").Value = 536
$wsAll.Range("F33").Value = 166
$wsAll.Range("F34").Value = 13271
$wsAll.Range("F39").Value = 5462
$wsAll.Range("F40").Value = 5565
$wsAll.Range("F41").Value = 35
